# "Spelling and grammar check from Word"
#
# 1) The "Block Text" example paragraph had been split into two runs by a
#    stray "_GoBack" bookmark sitting in the middle of the word "Text"
#    (" Block T" + bookmarkStart/_GoBack/bookmarkEnd + "ext. "). A
#    spelling/grammar pass re-merges that into a single run reading
#    " Block Text. " and drops the bookmark from that spot.
# 2) The "_GoBack" bookmark reappears at the new last edit location: a
#    new (empty) "Definition"-styled paragraph right after the existing
#    Definition/DefinitionTerm examples, followed by a new
#    "Bibliography"-styled paragraph containing the text
#    "Kumler and others, 2019".
# 3) The "Bibliography" style picks up <w:noProof/> (Word marking the
#    style as now in use / spell-check exempt).

$d = $word.ActiveDocument

# --- 1) Merge the "Block Text" runs back into one, dropping the stray
#        mid-word _GoBack bookmark. ---
$blockTextRange = $d.Content
$foundBlockText = $blockTextRange.Find.Execute(
    " Block Text. ", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)
if (-not $foundBlockText) {
    throw "Could not find the 'Block Text' example paragraph"
}
# Force a real content rewrite (same text in, same text out) so the
# engine collapses the two runs / removes the bookmark that splits them.
$blockTextRange.Text = "TEMP_BLOCK_TEXT_PLACEHOLDER"
$blockTextRange.Text = " Block Text. "

# --- 2) Find the last "Definition" example paragraph (there are two
#        Definition/DefinitionTerm pairs; we want the final one). ---
$defRange = $d.Content
$lastDefStart = -1
$lastDefEnd = -1
while ($defRange.Find.Execute(" Definition ", $true, $false, $false, $false,
                               $false, $true, 1, $false, "", 0)) {
    $lastDefStart = $defRange.Start
    $lastDefEnd = $defRange.End
    $defRange.Collapse(0)
}
if ($lastDefStart -eq -1) {
    throw "Could not find the 'Definition' example paragraph"
}

# Insert right after that paragraph's closing paragraph mark: a new
# empty Definition-styled paragraph carrying the _GoBack bookmark, then
# a new Bibliography-styled paragraph with the citation text.
$insertionPoint = $d.Range($lastDefEnd + 1, $lastDefEnd + 1)
$newParasXml = @'
<?xml version="1.0"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Definition"/>
            </w:pPr>
            <w:bookmarkStart w:id="12" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="12"/>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Bibliography"/>
            </w:pPr>
            <w:r>
              <w:t>Kumler and others, 2019</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@
[void]$insertionPoint.InsertXML($newParasXml)

# --- 3) Mark the Bibliography style as no-proof (now that it's in use). ---
$bibliographyStyle = $d.Styles("Bibliography")
$bibliographyStyle.NoProofing = $true

Write-Output "Done: merged Block Text runs, added _GoBack + bibliography paragraphs, updated Bibliography style."
